$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with new column titles
$ws.Range("A1").Value = " Lectura actual"
$ws.Range("B1").Value = " Lectura anterior"
$ws.Range("C1").Value = " Consumo (m3)"
$ws.Range("D1").Value = " Vertimiento (m3)"
$ws.Range("E1").Value = "Total a pagar"

# Update the single remaining data row (row 2) with the values that used
# to live in row 3, columns C:G
$ws.Range("A2").Value = 492
$ws.Range("B2").Value = 485
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 43500

# Remove everything else: extra columns F:G and extra rows 3:5
$ws.Range("F1:G5").Clear()
$ws.Range("A3:E5").Clear()
